$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 2 (A2:G2) already has values for the first data row.
# We need to rewrite row 2 and add rows 3-11 with the new data pattern.

$startValue = 12
for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $startValue - $i
    $ws.Cells.Item($row, 2).Value = 6
    $ws.Cells.Item($row, 3).Value = 776
    $ws.Cells.Item($row, 4).Value = 2
    $ws.Cells.Item($row, 5).Value = "С"
    $ws.Cells.Item($row, 6).Value = 0
    $ws.Cells.Item($row, 7).Value = 0
}
